# The workbook gained one new price record for "Pepino ensalada" at
# Macroferia Regional de Talca. The new record is inserted as row 138,
# pushing the former rows 138-215 down to rows 139-216 (dimension grows
# from A1:R215 to A1:R216).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 138 (shifts 138..215 down to 139..216)
$ws.Rows("138").Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A138").Value = 5
$ws.Range("B138").Value = "Macroferia Regional de Talca"
$ws.Range("C138").Value = "Maule"
$ws.Range("D138").Value = 44455
$ws.Range("E138").Value = 7
$ws.Range("F138").Value = 100112043
$ws.Range("G138").Value = "Pepino ensalada"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 300
$ws.Range("K138").Value = 16000
$ws.Range("L138").Value = 16000
$ws.Range("M138").Value = 16000
$ws.Range("N138").Value = '$/caja 60 unidades'
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 267
$ws.Range("Q138").Value = 60
$ws.Range("R138").Value = "Hortaliza"
